# Update NATMI LR-pair TPM output values (Sema4a-Plxnb1) to reflect
# recomputed statistics from the updated TPM scripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.361143666666667
$ws.Range("H2").Value = 28.083431
$ws.Range("I2").Value = 0.1965934861218525
$ws.Range("J2").Value = 0.1965934861218526
$ws.Range("M2").Value = 1.461859
$ws.Range("N2").Value = 4.385577
$ws.Range("O2").Value = 0.1790970628189019
$ws.Range("P2").Value = 0.1790970628189019
$ws.Range("Q2").Value = 13.68467211940967
$ws.Range("R2").Value = 123.162049074687
$ws.Range("S2").Value = 0.03520931593375234
$ws.Range("T2").Value = 0.03520931593375234
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.361143666666667
$ws.Range("H3").Value = 28.083431
$ws.Range("I3").Value = 0.1965934861218525
$ws.Range("J3").Value = 0.1965934861218526
$ws.Range("O3").Value = 0.1506444277390854
$ws.Range("P3").Value = 0.1506444277390854
$ws.Range("Q3").Value = 11.51062763273811
$ws.Range("R3").Value = 103.595648694643
$ws.Range("S3").Value = 0.0296157132140583
$ws.Range("T3").Value = 0.0296157132140583
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.361143666666667
$ws.Range("H4").Value = 28.083431
$ws.Range("I4").Value = 0.1965934861218525
$ws.Range("J4").Value = 0.1965934861218526
$ws.Range("M4").Value = 5.443148666666667
$ws.Range("N4").Value = 16.329446
$ws.Range("O4").Value = 0.6668577056245659
$ws.Range("P4").Value = 0.6668577056245659
$ws.Range("Q4").Value = 50.95409666769178
$ws.Range("R4").Value = 458.5868700092261
$ws.Range("S4").Value = 0.1310998810959535
$ws.Range("T4").Value = 0.1310998810959536
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.361143666666667
$ws.Range("H5").Value = 28.083431
$ws.Range("I5").Value = 0.1965934861218525
$ws.Range("J5").Value = 0.1965934861218526
$ws.Range("M5").Value = 0.02775866666666667
$ws.Range("N5").Value = 0.083276
$ws.Range("O5").Value = 0.003400803817446799
$ws.Range("P5").Value = 0.003400803817446798
$ws.Range("Q5").Value = 0.2598528666617778
$ws.Range("R5").Value = 2.338675799956
$ws.Range("S5").Value = 0.0006685758780883704
$ws.Range("T5").Value = 0.0006685758780883704
$ws.Range("I6").Value = 0.2809659460057553
$ws.Range("J6").Value = 0.2809659460057554
$ws.Range("M6").Value = 1.461859
$ws.Range("N6").Value = 4.385577
$ws.Range("O6").Value = 0.1790970628189019
$ws.Range("P6").Value = 0.1790970628189019
$ws.Range("Q6").Value = 19.55775302456033
$ws.Range("R6").Value = 176.019777221043
$ws.Range("S6").Value = 0.05032017568176497
$ws.Range("T6").Value = 0.05032017568176497
$ws.Range("I7").Value = 0.2809659460057553
$ws.Range("J7").Value = 0.2809659460057554
$ws.Range("O7").Value = 0.1506444277390854
$ws.Range("P7").Value = 0.1506444277390854
$ws.Range("S7").Value = 0.04232595415020778
$ws.Range("T7").Value = 0.04232595415020778
$ws.Range("I8").Value = 0.2809659460057553
$ws.Range("J8").Value = 0.2809659460057554
$ws.Range("M8").Value = 5.443148666666667
$ws.Range("N8").Value = 16.329446
$ws.Range("O8").Value = 0.6668577056245659
$ws.Range("P8").Value = 0.6668577056245659
$ws.Range("Q8").Value = 72.8221786770349
$ws.Range("R8").Value = 655.3996080933141
$ws.Range("S8").Value = 0.1873643061120337
$ws.Range("T8").Value = 0.1873643061120337
$ws.Range("I9").Value = 0.2809659460057553
$ws.Range("J9").Value = 0.2809659460057554
$ws.Range("M9").Value = 0.02775866666666667
$ws.Range("N9").Value = 0.083276
$ws.Range("O9").Value = 0.003400803817446799
$ws.Range("P9").Value = 0.003400803817446798
$ws.Range("Q9").Value = 0.371374494364889
$ws.Range("R9").Value = 3.342370449284001
$ws.Range("S9").Value = 0.0009555100617489237
$ws.Range("T9").Value = 0.0009555100617489238
$ws.Range("G10").Value = 7.684952333333334
$ws.Range("H10").Value = 23.054857
$ws.Range("I10").Value = 0.1613917725961189
$ws.Range("J10").Value = 0.1613917725961189
$ws.Range("M10").Value = 1.461859
$ws.Range("N10").Value = 4.385577
$ws.Range("O10").Value = 0.1790970628189019
$ws.Range("P10").Value = 0.1790970628189019
$ws.Range("Q10").Value = 11.23431673305433
$ws.Range("R10").Value = 101.108850597489
$ws.Range("S10").Value = 0.02890479243510104
$ws.Range("T10").Value = 0.02890479243510103
$ws.Range("G11").Value = 7.684952333333334
$ws.Range("H11").Value = 23.054857
$ws.Range("I11").Value = 0.1613917725961189
$ws.Range("J11").Value = 0.1613917725961189
$ws.Range("O11").Value = 0.1506444277390854
$ws.Range("P11").Value = 0.1506444277390854
$ws.Range("Q11").Value = 9.44955315655789
$ws.Range("R11").Value = 85.045978409021
$ws.Range("S11").Value = 0.02431277122453894
$ws.Range("T11").Value = 0.02431277122453893
$ws.Range("G12").Value = 7.684952333333334
$ws.Range("H12").Value = 23.054857
$ws.Range("I12").Value = 0.1613917725961189
$ws.Range("J12").Value = 0.1613917725961189
$ws.Range("M12").Value = 5.443148666666667
$ws.Range("N12").Value = 16.329446
$ws.Range("O12").Value = 0.6668577056245659
$ws.Range("P12").Value = 0.6668577056245659
$ws.Range("Q12").Value = 41.83033804658023
$ws.Range("R12").Value = 376.4730424192221
$ws.Range("S12").Value = 0.1076253471801295
$ws.Range("T12").Value = 0.1076253471801295
$ws.Range("G13").Value = 7.684952333333334
$ws.Range("H13").Value = 23.054857
$ws.Range("I13").Value = 0.1613917725961189
$ws.Range("J13").Value = 0.1613917725961189
$ws.Range("M13").Value = 0.02775866666666667
$ws.Range("N13").Value = 0.083276
$ws.Range("O13").Value = 0.003400803817446799
$ws.Range("P13").Value = 0.003400803817446798
$ws.Range("Q13").Value = 0.2133240301702223
$ws.Range("R13").Value = 1.919916271532
$ws.Range("S13").Value = 0.0005488617563493868
$ws.Range("T13").Value = 0.0005488617563493866
$ws.Range("G14").Value = 17.19197166666666
$ws.Range("H14").Value = 51.57591499999999
$ws.Range("I14").Value = 0.3610487952762732
$ws.Range("J14").Value = 0.3610487952762733
$ws.Range("M14").Value = 1.461859
$ws.Range("N14").Value = 4.385577
$ws.Range("O14").Value = 0.1790970628189019
$ws.Range("P14").Value = 0.1790970628189019
$ws.Range("Q14").Value = 25.13223850866166
$ws.Range("R14").Value = 226.1901465779549
$ws.Range("S14").Value = 0.06466277876828357
$ws.Range("T14").Value = 0.06466277876828357
$ws.Range("G15").Value = 17.19197166666666
$ws.Range("H15").Value = 51.57591499999999
$ws.Range("I15").Value = 0.3610487952762732
$ws.Range("J15").Value = 0.3610487952762733
$ws.Range("O15").Value = 0.1506444277390854
$ws.Range("P15").Value = 0.1506444277390854
$ws.Range("Q15").Value = 21.13955208616611
$ws.Range("R15").Value = 190.255968775495
$ws.Range("S15").Value = 0.05438998915028038
$ws.Range("T15").Value = 0.05438998915028038
$ws.Range("G16").Value = 17.19197166666666
$ws.Range("H16").Value = 51.57591499999999
$ws.Range("I16").Value = 0.3610487952762732
$ws.Range("J16").Value = 0.3610487952762733
$ws.Range("M16").Value = 5.443148666666667
$ws.Range("N16").Value = 16.329446
$ws.Range("O16").Value = 0.6668577056245659
$ws.Range("P16").Value = 0.6668577056245659
$ws.Range("Q16").Value = 93.57845765478775
$ws.Range("R16").Value = 842.20611889309
$ws.Range("S16").Value = 0.2407681712364492
$ws.Range("T16").Value = 0.2407681712364492
$ws.Range("G17").Value = 17.19197166666666
$ws.Range("H17").Value = 51.57591499999999
$ws.Range("I17").Value = 0.3610487952762732
$ws.Range("J17").Value = 0.3610487952762733
$ws.Range("M17").Value = 0.02775866666666667
$ws.Range("N17").Value = 0.083276
$ws.Range("O17").Value = 0.003400803817446799
$ws.Range("P17").Value = 0.003400803817446798
$ws.Range("Q17").Value = 0.4772262108377777
$ws.Range("R17").Value = 4.29503589754
$ws.Range("S17").Value = 0.001227856121260118
$ws.Range("T17").Value = 0.001227856121260118
